# "Add item" — appends a new inventory row (row 5) to the Inventory sheet,
# mirroring the existing data rows 2-4 (Item Desc, Category, Made in,
# Size/ml-g-oz, Unit, Quantity, Bar Code, Unit Price, Profit, Unit Price
# after profit, VAT, Total, On hand qty, Sold qty, Total sales).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last data row's formatting down into the new row 5 first
# (copy row 4, insert the copy at row 5) so the new row picks up the same
# cell style the rest of the table uses, then overwrite the copied values
# below with the real new-item data.
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(5).Insert()

$row = 5

$ws.Cells.Item($row, 1).Value  = "emp"             # Item Desc
$ws.Cells.Item($row, 2).Value  = "stock"            # Category
$ws.Cells.Item($row, 3).Value  = "chn"              # Made in
$ws.Cells.Item($row, 4).Value  = 10                 # Size/ml-g-oz
$ws.Cells.Item($row, 5).Value  = "pcs"               # Unit
$ws.Cells.Item($row, 6).Value  = 249                # Quantity
$ws.Cells.Item($row, 7).Value  = 951235648952       # Bar Code
$ws.Cells.Item($row, 8).Value  = 12.35              # Unit Price

# These columns hold numeric-looking text in the source data (stored as
# strings, not numbers) — a leading apostrophe forces text storage while
# keeping the cell's General number format, matching the other rows.
$ws.Cells.Item($row, 9).Value  = "'5.56"            # Profit
$ws.Cells.Item($row, 10).Value = "'17.91"           # Unit Price after profit
$ws.Cells.Item($row, 11).Value = "'2.69"            # VAT
$ws.Cells.Item($row, 12).Value = "'20.59"           # Total

$ws.Cells.Item($row, 13).Value = 249                # On hand qty
$ws.Cells.Item($row, 14).Value = 0                  # Sold qty
$ws.Cells.Item($row, 15).Value = 0                  # Total sales
